# Apply updated Cv (col B) and Xt (col D) values on the "Valve_2.0_600_1"
# sheet, then move the active selection to G10 (single cell), matching
# the author's re-upload of the valve database.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valve_2.0_600_1")

# New Cv values (column B, rows 2-12)
$bValues = @(0, 20, 30, 50, 90, 120, 150, 180, 210, 300, 320)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# New Xt values (column D, rows 2-12) - all become 0.688
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 4).Value = 0.688
}

# Update the active sheet's selection to G10
$ws.Activate()
$ws.Range("G10").Select()
